# Apply cryptocurrency price/volume updates to Sheet1
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '70.786.88'
$ws.Range('E2').Value = '  +1.93%  '
$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '3.582.63'
$ws.Range('E3').Value = '  +1.17%  '
$ws.Range('D4').NumberFormat = "@"
$ws.Range('D4').Value = '0.999'
$ws.Range('E4').Value = '  -0.03%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '607.20'
$ws.Range('E5').Value = '  +4.48%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '174.96'
$ws.Range('E6').Value = '  +1.22%  '
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '3.577.78'
$ws.Range('E7').Value = '  +1.24%  '
$ws.Range('E8').Value = '  +0.41%  '
$ws.Range('E9').Value = '  -0.04%  '
$ws.Range('E10').Value = '  +3.56%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '7.48'
$ws.Range('E11').Value = '  +10.35%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '0.590'
$ws.Range('E12').Value = '  +0.61%  '
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '47.13'
$ws.Range('E13').Value = '  -0.66%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '0.0000279'
$ws.Range('E14').Value = '  +1.00%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '4.160.64'
$ws.Range('E15').Value = '  +1.26%  '
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '8.45'
$ws.Range('E16').Value = '  -1.34%  '
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '619.17'
$ws.Range('E17').Value = '  -1.69%  '
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '3.566.07'
$ws.Range('E18').Value = '  +0.96%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '70.868.08'
$ws.Range('E19').Value = '  +2.02%  '
$ws.Range('E20').Value = '  -2.36%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '17.52'
$ws.Range('E21').Value = '  +0.63%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '0.891'
$ws.Range('E22').Value = '  -0.39%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '9.40'
$ws.Range('E23').Value = '  -16.26%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '16.23'
$ws.Range('E24').Value = '  +1.26%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '97.88'
$ws.Range('E25').Value = '  +0.17%  '
$ws.Range('E26').Value = '  +0.03%  '
$ws.Range('E27').Value = '  +0.05%  '
$ws.Range('E28').Value = '  +0.39%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '9.33'
$ws.Range('E29').Value = '  -0.23%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '33.63'
$ws.Range('E30').Value = '  +1.89%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '8.54'
$ws.Range('E31').Value = '  -0.57%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '3.07'
$ws.Range('E32').Value = '  -3.13%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '7.11'
$ws.Range('E33').Value = '  +1.15%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '1.31'
$ws.Range('E34').Value = '  -2.56%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '3.81'
$ws.Range('E35').Value = '  +8.14%  '
$ws.Range('E36').Value = '  -0.71%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '10.90'
$ws.Range('E37').Value = '  +0.74%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '0.0490'
$ws.Range('E38').Value = '  +6.18%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '57.42'
$ws.Range('E39').Value = '  +0.04%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '0.999'
$ws.Range('E40').Value = '  +0.00%  '
$ws.Range('E41').Value = '  +3.84%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '3.411.55'
$ws.Range('E42').Value = '  +0.29%  '
$ws.Range('E43').Value = '  -1.71%  '
$ws.Range('B44').Value = 'PEPE'
$ws.Range('C44').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '0.0₃0715'
$ws.Range('E44').Value = '  +1.50%  '
$ws.Range('B45').Value = 'ThetaToken'
$ws.Range('C45').Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '3.01'
$ws.Range('E45').Value = '  +8.86%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '33.15'
$ws.Range('E46').Value = '  +0.95%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '2.69'
$ws.Range('E47').Value = '  +4.72%  '
$ws.Range('E48').Value = '  +0.57%  '
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '133.15'
$ws.Range('E49').Value = '  -0.02%  '
$ws.Range('E51').Value = '  +0.90%  '
